$d = $word.ActiveDocument

# --- Change 1 ------------------------------------------------------------
# "TODOS LOS PAGOS DEBERÁN DE SER REALIZADOS A BENEFICIO “{{SEXO_8}} PROMITENTE {{SEXO_2}}”"
# loses the opening curly quote and the space before "PROMITENTE" (closing quote stays).
$oldText1 = "TODOS LOS PAGOS DEBERÁN DE SER REALIZADOS A BENEFICIO " + [char]8220 + "{{SEXO_8}} PROMITENTE {{SEXO_2}}"
$newText1 = "TODOS LOS PAGOS DEBERÁN DE SER REALIZADOS A BENEFICIO {{SEXO_8}}PROMITENTE {{SEXO_2}}"
$d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)

# --- Change 2 ------------------------------------------------------------
# Payment schedule sentence rewritten (days/hours updated).
$oldText2 = "LOS PAGOS DEBERÁN REALIZARSE DE LUNES A SÁBADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. PARA EFECTUAR UN PAGO EN DOMINGO, SERÁ INDISPENSABLE PROGRAMAR UNA CITA CON TRES DÍAS DE ANTICIPACIÓN. CADA PAGO DEBERÁ SER NOTIFICADO Y CONFIRMADO AL NÚMERO TELEFÓNICO 951 189 9298."
$newText2 = "LOS PAGOS DEBERÁN REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., Y LOS SÁBADOS DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SERÁ INDISPENSABLE PROGRAMAR UNA CITA CON AL MENOS TRES DÍAS DE ANTICIPACIÓN. CADA PAGO DEBERÁ SER NOTIFICADO Y CONFIRMADO AL NÚMERO TELEFÓNICO 951 189 9298."
$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)

# --- Change 3 ------------------------------------------------------------
# Append ' {{SEXO_7}}PROMITENTE {{SEXO_2}}”' right after the penalty-clause
# bookmark, before the closing period of that sentence.
$rng3 = $d.Content
$rng3.Find.Execute("CUBRIR LA PENA CONVENCIONAL ESTABLECIDA POR INCUMPLIMIENTO, SIN PERJUICIO DE OTRAS ACCIONES LEGALES QUE EN SU CASO CORRESPONDAN")
$rng3.Collapse(0)
$rng3.InsertAfter(" {{SEXO_7}}PROMITENTE {{SEXO_2}}" + [char]8221)

Write-Host "Edits applied"
